# Daily KHL injuries refresh:
#  - "snapshot": drop rows for players who returned to play, insert rows for
#    newly injured players (kept in team/player sorted order), and refresh the
#    scraped_at timestamp on every remaining row.
#  - "returned": append an audit row per player who left the snapshot.
#  - "new_injured": append an audit row per player newly added to the snapshot.

$wb = $excel.ActiveWorkbook
$snapshot = $wb.Worksheets.Item("snapshot")
$returned = $wb.Worksheets.Item("returned")
$newInjured = $wb.Worksheets.Item("new_injured")

# --- 1) Remove rows for players who are no longer injured -------------------
$snapshot.Rows.Item(45).Delete() | Out-Null
$snapshot.Rows.Item(43).Delete() | Out-Null
$snapshot.Rows.Item(23).Delete() | Out-Null

# --- 2) Insert rows for newly injured players (alphabetical club/player order)
$snapshot.Rows.Item(19).Insert() | Out-Null
$newRow = $snapshot.Range("A19:K19")
$newRow.NumberFormat = "@"
$snapshot.Cells.Item(19, 1).Value = 'ЛОК'
$snapshot.Cells.Item(19, 2).Value = 'Локомотив'
$snapshot.Cells.Item(19, 3).Value = 'lokomotiv'
$snapshot.Cells.Item(19, 4).Value = 'Паник Рихард'
$snapshot.Cells.Item(19, 5).Value = '14'
$snapshot.Cells.Item(19, 6).Value = 'нападающий'
$snapshot.Cells.Item(19, 7).Value = '16071'
$snapshot.Cells.Item(19, 8).Value = '1369_ЛОК_паникрихард'
$snapshot.Cells.Item(19, 9).Value = 'injured_active'
$snapshot.Cells.Item(19, 10).Value = 'https://www.khl.ru/clubs/lokomotiv/team/'
$snapshot.Cells.Item(19, 11).Value = '2025-11-08T03:03:14.245871+00:00'

# --- 3) Refresh scraped_at (column K) for every row, matched by player_uid ---
$scrapedAt = @{
    '1369_АВТ_зборовскийсергей' = '2025-11-08T03:02:53.664217+00:00'
    '1369_АВТ_кизимовсемен' = '2025-11-08T03:02:53.664252+00:00'
    '1369_АВТ_трямкинникита' = '2025-11-08T03:02:53.664277+00:00'
    '1369_АДМ_грманмарио' = '2025-11-08T03:02:56.438192+00:00'
    '1369_АДМ_старковстепан' = '2025-11-08T03:02:56.438221+00:00'
    '1369_АДМ_шепелевалександр' = '2025-11-08T03:02:56.438240+00:00'
    '1369_АКБ_яруллинальберт' = '2025-11-08T03:02:59.346418+00:00'
    '1369_АМР_абросимовроман' = '2025-11-08T03:03:01.737685+00:00'
    '1369_АМР_броадхерсталекс' = '2025-11-08T03:03:01.737715+00:00'
    '1369_АМР_гиздатуллинартур' = '2025-11-08T03:03:01.737733+00:00'
    '1369_БАР_бояркинникита' = '2025-11-08T03:03:04.537647+00:00'
    '1369_БАР_галимовэмиль' = '2025-11-08T03:03:04.537675+00:00'
    '1369_БАР_мухаметовмаксим' = '2025-11-08T03:03:04.537692+00:00'
    '1369_БАР_уотерспунтайлер' = '2025-11-08T03:03:04.537708+00:00'
    '1369_ДМН_уэллексавье' = '2025-11-08T03:03:09.141509+00:00'
    '1369_ЛАД_ожгихиналексей' = '2025-11-08T03:03:11.477838+00:00'
    '1369_ЛОК_волковалександрс' = '2025-11-08T03:03:14.245804+00:00'
    '1369_ЛОК_сергеевандрей' = '2025-11-08T03:03:14.245899+00:00'
    '1369_ММГ_козловандрейе' = '2025-11-08T03:03:16.521996+00:00'
    '1369_НХК_дергачевалександр' = '2025-11-08T03:03:19.232236+00:00'
    '1369_НХК_попугаевникитао' = '2025-11-08T03:03:19.232266+00:00'
    '1369_СЕВ_ващенкогригорий' = '2025-11-08T03:03:22.151013+00:00'
    '1369_СЕВ_грудининвладимир' = '2025-11-08T03:03:22.151047+00:00'
    '1369_СЕВ_цицюравладислав' = '2025-11-08T03:03:22.151067+00:00'
    '1369_СИБ_гордеевфедор' = '2025-11-08T03:03:25.118685+00:00'
    '1369_СИБ_калиниченкороман' = '2025-11-08T03:03:25.118718+00:00'
    '1369_СИБ_прискичейзэванс' = '2025-11-08T03:03:25.118744+00:00'
    '1369_СИБ_пьяноввалентин' = '2025-11-08T03:03:25.118762+00:00'
    '1369_СИБ_широковсергей' = '2025-11-08T03:03:25.118782+00:00'
    '1369_СКА_зайцевникитаи' = '2025-11-08T03:03:27.873012+00:00'
    '1369_СОЧ_гуськовматвей' = '2025-11-08T03:03:30.152503+00:00'
    '1369_СОЧ_мачулинвасилий' = '2025-11-08T03:03:30.152533+00:00'
    '1369_СОЧ_хомченкопавел' = '2025-11-08T03:03:30.152552+00:00'
    '1369_СПР_воробьевиванв' = '2025-11-08T03:03:32.406671+00:00'
    '1369_СПР_порядинпавел' = '2025-11-08T03:03:32.406699+00:00'
    '1369_СПР_рубцовгерман' = '2025-11-08T03:03:32.406716+00:00'
    '1369_СЮЛ_алалыкинданил' = '2025-11-08T03:03:35.218024+00:00'
    '1369_СЮЛ_берлевантон' = '2025-11-08T03:03:35.218053+00:00'
    '1369_СЮЛ_зоркинникита' = '2025-11-08T03:03:35.218073+00:00'
    '1369_СЮЛ_кузьминглеб' = '2025-11-08T03:03:35.218091+00:00'
    '1369_СЮЛ_хворовниколай' = '2025-11-08T03:03:35.218107+00:00'
    '1369_СЮЛ_янденис' = '2025-11-08T03:03:35.218121+00:00'
    '1369_ТОР_науменковмихаил' = '2025-11-08T03:03:38.123971+00:00'
    '1369_ТОР_рожковникитаа' = '2025-11-08T03:03:38.124004+00:00'
    '1369_ЦСК_моисеевданила' = '2025-11-08T03:03:43.722175+00:00'
    '1369_ЦСК_уильямсколби' = '2025-11-08T03:03:43.722204+00:00'
    '1369_ШДР_гроложереми' = '2025-11-08T03:03:46.028334+00:00'
    '1369_ШДР_саттеррайли' = '2025-11-08T03:03:46.028363+00:00'
}

$lastRow = $snapshot.Cells.Item($snapshot.Rows.Count, 8).End(-4162).Row  # xlUp
for ($r = 2; $r -le $lastRow; $r++) {
    $uid = $snapshot.Cells.Item($r, 8).Text
    if ($scrapedAt.ContainsKey($uid)) {
        $snapshot.Cells.Item($r, 11).Value = $scrapedAt[$uid]
    }
}

# --- 4) Log returned players ---------------------------------------------------
$returnedRow = $returned.Cells.Item($returned.Rows.Count, 1).End(-4162).Row + 1
if ($returned.Cells.Item(1,1).Text -eq "") { $returnedRow = 2 }
$rowRange = $returned.Range("A" + $returnedRow + ":G" + $returnedRow)
$rowRange.NumberFormat = "@"
$returned.Cells.Item($returnedRow, 1).Value = 'НХК'
$returned.Cells.Item($returnedRow, 2).Value = 'Нефтехимик'
$returned.Cells.Item($returnedRow, 3).Value = 'Профака Лука'
$returned.Cells.Item($returnedRow, 4).Value = '1369_НХК_профакалука'
$returned.Cells.Item($returnedRow, 5).Value = 'RETURN'
$returned.Cells.Item($returnedRow, 6).Value = '2025-11-08T11:03:46.545097+08:00'
$returned.Cells.Item($returnedRow, 7).Value = '2025-11-08'
$returnedRow++

$rowRange = $returned.Range("A" + $returnedRow + ":G" + $returnedRow)
$rowRange.NumberFormat = "@"
$returned.Cells.Item($returnedRow, 1).Value = 'СЮЛ'
$returned.Cells.Item($returnedRow, 2).Value = 'Салават Юлаев'
$returned.Cells.Item($returnedRow, 3).Value = 'Пименов Артём'
$returned.Cells.Item($returnedRow, 4).Value = '1369_СЮЛ_пименовартем'
$returned.Cells.Item($returnedRow, 5).Value = 'RETURN'
$returned.Cells.Item($returnedRow, 6).Value = '2025-11-08T11:03:46.545097+08:00'
$returned.Cells.Item($returnedRow, 7).Value = '2025-11-08'
$returnedRow++

$rowRange = $returned.Range("A" + $returnedRow + ":G" + $returnedRow)
$rowRange.NumberFormat = "@"
$returned.Cells.Item($returnedRow, 1).Value = 'СЮЛ'
$returned.Cells.Item($returnedRow, 2).Value = 'Салават Юлаев'
$returned.Cells.Item($returnedRow, 3).Value = 'Хохряков Пётр'
$returned.Cells.Item($returnedRow, 4).Value = '1369_СЮЛ_хохряковпетр'
$returned.Cells.Item($returnedRow, 5).Value = 'RETURN'
$returned.Cells.Item($returnedRow, 6).Value = '2025-11-08T11:03:46.545097+08:00'
$returned.Cells.Item($returnedRow, 7).Value = '2025-11-08'
$returnedRow++

# --- 5) Log newly injured players -----------------------------------------------
$newInjuredRow = $newInjured.Cells.Item($newInjured.Rows.Count, 1).End(-4162).Row + 1
if ($newInjured.Cells.Item(1,1).Text -eq "") { $newInjuredRow = 2 }
$rowRange2 = $newInjured.Range("A" + $newInjuredRow + ":G" + $newInjuredRow)
$rowRange2.NumberFormat = "@"
$newInjured.Cells.Item($newInjuredRow, 1).Value = 'ЛОК'
$newInjured.Cells.Item($newInjuredRow, 2).Value = 'Локомотив'
$newInjured.Cells.Item($newInjuredRow, 3).Value = 'Паник Рихард'
$newInjured.Cells.Item($newInjuredRow, 4).Value = '1369_ЛОК_паникрихард'
$newInjured.Cells.Item($newInjuredRow, 5).Value = 'INJURED_NEW'
$newInjured.Cells.Item($newInjuredRow, 6).Value = '2025-11-08T11:03:46.545097+08:00'
$newInjured.Cells.Item($newInjuredRow, 7).Value = '2025-11-08'
$newInjuredRow++

